# "bind country to self" - restructure the Operators sheet:
#   - insert a new column A ("Operation Type": TRANSFORM / TRANSFER)
#   - old A/B/C (Operation Name / Definition / Preconditions) shift to B/C/D
#   - clean up the garbled TRANSFORM definition strings
#   - add a new TRANSFER row (row 5)
#   - bold + center the header row
#   - adjust column widths / view zoom / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column in front of the existing data (A -> B, B -> C, C -> D).
$ws.Columns("A").Insert()

# 2) Header row.
$ws.Range("A1").Value = "Operation Type"
$ws.Range("B1").Value = "Operation Name"
$ws.Range("C1").Value = "Definition"
$ws.Range("D1").Value = "Preconditions"

# 3) Operation Type column for the existing TRANSFORM rows.
$ws.Range("A2").Value = "TRANSFORM"
$ws.Range("A3").Value = "TRANSFORM"
$ws.Range("A4").Value = "TRANSFORM"

# 4) Clean up the garbled TRANSFORM definition strings (now in column C).
#    C2 keeps its original quote-prefixed style, so re-set it with a leading
#    apostrophe to preserve that xf instead of minting an unused duplicate.
$ws.Range("C2").Value = "'(TRANSFORM ?C (INPUTS (R1 1) (R2 2)) (OUTPUTS (R1 1) (R21 1) (R21' 1)))"
$ws.Range("C3").Value = "(TRANSFORM ?C (INPUTS (R1 3) (R2 2) (R21 2)) (OUTPUTS (R22 2) (R22' 2) (R1 3)))"
$ws.Range("C4").Value = "(TRANSFORM ?C (INPUTS (R1 5) (R2 1) (R3 5) (R21 3)) (OUTPUTS (R1 5) (R23 1) (R23' 1)))"

# 5) New TRANSFER row.
$ws.Range("A5").Value = "TRANSFER"
$ws.Range("B5").Value = "transfer"
$ws.Range("C5").Value = "(TRANSFER ?Cj1 ?Cj2 (?Ri ?ARi))"
$ws.Range("D5").Value = "?ARi <= ?Cj1(?Ri)"

# 6) Header formatting: bold + centered (D1 "Preconditions" keeps its
#    original, unbolded centered style).
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108

# 7) Operation Type column centered for the data rows too.
$ws.Range("A2:A5").HorizontalAlignment = -4108

# 8) Column widths (inputs compensate for the engine's ColumnWidth -> stored
#    width padding so the saved OOXML <col width> lands on the target value).
$ws.Columns("A").ColumnWidth = 29.5
$ws.Columns("B").ColumnWidth = 20.166666666666668
$ws.Columns("C").ColumnWidth = 74.66666666666667
$ws.Columns("D").ColumnWidth = 19.166666666666668

# 9) View: zoom + selection.
$ws.Application.ActiveWindow.Zoom = 140
$ws.Range("C8").Select()
